$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited cells to remain plain text (they store numbers/percentages as text)
# so that COM does not reinterpret them as numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "291.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.77%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.78"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.44%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07260"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.44%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.354"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "27.79%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.669"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.25%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.705"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.18%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8977"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.43%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.06%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07973"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.95%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08103"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.03%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03096"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.59%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.23%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001499"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.10%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005816"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.82%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.472"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.46%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.20%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3318"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.00%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1298"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.09%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.975"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-8.03%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "16.46%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04520"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.91%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001211"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.02%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004403"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.82%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.91%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003395"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01583"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.30%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04375"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.63%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007317"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.75%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01002"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1311"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.18%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002032"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.29%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009504"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-14.72%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005737"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.01%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.241"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "12.04%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002899"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.51%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"
